# Insert a new weekly/monthly price record at row 75 ("Ajo" / "Chino" / "Primera",
# Vega Monumental Concepción), pushing the existing rows 75-108 down to 76-109.
# The former row 108 therefore becomes the new row 109 without any further edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 75..108 down to 76..109 by inserting a new blank row at 75.
$ws.Rows.Item(75).EntireRow.Insert()

# Populate the newly inserted row 75 with the new record's data.
$ws.Range("A75").Value = 11
$ws.Range("B75").Value = "Vega Monumental Concepción"
$ws.Range("C75").Value = "Bíobío"
$ws.Range("D75").Value = 44523
$ws.Range("E75").Value = 8
$ws.Range("F75").Value = 100112003
$ws.Range("G75").Value = "Ajo"
$ws.Range("H75").Value = "Chino"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 400
$ws.Range("K75").Value = 18000
$ws.Range("L75").Value = 19000
$ws.Range("M75").Value = 18500
$ws.Range("N75").Value = "$/caja 10 kilos"
$ws.Range("O75").Value = "China"
$ws.Range("P75").Value = 1850
$ws.Range("Q75").Value = 10
$ws.Range("R75").Value = "Hortaliza"
